$wb = $excel.ActiveWorkbook

# --- Update the view state of the existing "covmatrix" sheet -------------
$ws2 = $wb.Worksheets.Item("covmatrix")
$ws2.Activate() | Out-Null
$ws2.Range("C11").Select() | Out-Null

# --- Add the new "Sheet3" worksheet after the last existing sheet --------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Fill in the note/label grid column by column (G, then H, then I) so the
# shared-string table is built up in the same order as the source edit.
$ws3.Range("G10").Value = "身高"
$ws3.Range("G11").Value = "体重"
$ws3.Range("G12").Value = "性别代码(0,1)"

$ws3.Range("H10").Value = "年龄"
$ws3.Range("H11").Value = "最近一次英语成绩(百分制)"
$ws3.Range("H12").Value = "最近一次数学成绩(百分制)"

$ws3.Range("I10").Value = "视力"
$ws3.Range("I11").Value = "一周上网的时间(小时制)"
$ws3.Range("I12").Value = "亲密朋友个数"

# Bold font for the whole label grid (new cellXfs entry)
$ws3.Range("G10:I12").Font.Bold = $true

# Size the columns to fit their (wide, Chinese-text) contents
$ws3.Columns.Item(7).ColumnWidth = 12.83203125 - (5/6)
$ws3.Columns.Item(8).ColumnWidth = 24.83203125 - (5/6)
$ws3.Columns.Item(9).ColumnWidth = 22.83203125 - (5/6)

# Make Sheet3 the active sheet/tab and set its view state
$ws3.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 140
$ws3.Range("G6").Select() | Out-Null
